$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId=1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 400
$ws1.Range("F3").Value = 403
$ws1.Range("F4").Value = 2662
$ws1.Range("F5").Value = 1305
$ws1.Range("F8").Value = 558
$ws1.Range("F13").Value = 11009
$ws1.Range("F14").Value = 6296
$ws1.Range("F22").Value = 34
$ws1.Range("F23").Value = 221
$ws1.Range("F26").Value = 54
$ws1.Range("F35").Value = 1192
$ws1.Range("F36").Value = 187
$ws1.Range("F37").Value = 353
$ws1.Range("F38").Value = 146

# Sheet "演出" (sheetId=2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 3636

# Sheet "全部类型" (sheetId=4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 403
$ws4.Range("F7").Value = 2662
$ws4.Range("F12").Value = 1305
$ws4.Range("F14").Value = 558
$ws4.Range("F18").Value = 11009
$ws4.Range("F19").Value = 3636
$ws4.Range("F27").Value = 34
$ws4.Range("F28").Value = 221
$ws4.Range("F31").Value = 54
$ws4.Range("F42").Value = 1192
$ws4.Range("F44").Value = 187
$ws4.Range("F45").Value = 146
